$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EJ45")

$ws.Range("B2").Value = 23.61449426666159
$ws.Range("C2").Value = 0.499016247702294
$ws.Range("D2").Value = 48.22898853332319
$ws.Range("E2").Value = 0.001899645122801
$ws.Range("F2").Value = 8.328412141714798

$ws.Range("B3").Value = 23.84676650578297
$ws.Range("C3").Value = 0.669735083050551
$ws.Range("D3").Value = 48.69353301156595
$ws.Range("E3").Value = 0.002195087736976
$ws.Range("F3").Value = 4.822091633683793

$ws.Range("B4").Value = 28.8166289008596
$ws.Range("C4").Value = 0.641933963400083
$ws.Range("D4").Value = 58.6332578017192
$ws.Range("E4").Value = 0.003014653675459
$ws.Range("F4").Value = 3.865008219927484

$ws.Range("B5").Value = 22.52306418688367
$ws.Range("C5").Value = 0.770141832396215
$ws.Range("D5").Value = 46.04612837376735
$ws.Range("E5").Value = 0.001922449438888
$ws.Range("F5").Value = 7.560878692913667

$ws.Range("B6").Value = 23.86784155497488
$ws.Range("C6").Value = 0.722599681918588
$ws.Range("D6").Value = 48.73568310994975
$ws.Range("E6").Value = 0.002386997377656
$ws.Range("F6").Value = 5.280232667329749

$ws.Range("B7").Value = 28.6108312921836
$ws.Range("C7").Value = 0.367454292394605
$ws.Range("D7").Value = 58.22166258436721
$ws.Range("E7").Value = 0.003205318024707
$ws.Range("F7").Value = 3.90131933195202
